# Automatic update of files.
# - Column C ("Förändrad") bumps from 46059 to 46060 for every data row (2-13).
# - Rows 6-10 get their Beteckning/Datum/Area values re-ordered to reflect a
#   refreshed data pull (the underlying case ids stay the same, only which
#   row they land on changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump "Förändrad" (column C) for every data row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# Capture current (pre-shuffle) values for rows 6-10 so the row-reorder below
# reads consistently regardless of write order. Use Value2 to read (Value
# read-back is unreliable in this host), Value to write.
$rows = 6..10
$colA = @{}
$colB = @{}
$colG = @{}
foreach ($r in $rows) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

# New row -> source row mapping observed in the update.
$mapping = @{ 6 = 10; 7 = 8; 8 = 9; 9 = 7; 10 = 6 }

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $ws.Cells.Item($dest, 1).Value = $colA[$src]
    $ws.Cells.Item($dest, 2).Value = $colB[$src]
    $ws.Cells.Item($dest, 7).Value = $colG[$src]
}
